$wb = $excel.ActiveWorkbook

# Update status text from "Ready for handoff" to "In Translation"
# This shared string is used on:
#  - Overview sheet (zh-cn, de-de columns)
#  - zh-cn sheet (Status column)
#  - de-de sheet (Status column)

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# Shrink the affected "status" columns' widths to reflect the new, shorter text
# (mirrors what Excel's AutoFit would do after the text shrank from "Ready for
# handoff" to "In Translation"). Target raw OOXML width is ~13.41; ColumnWidth
# 12.5 maps to raw width 13.333333333333334 (raw = ColumnWidth +
# 0.8333333333333334), the closest value this engine can express.
$targetColumnWidth = 12.5
$wsOverview.Columns.Item(5).ColumnWidth = $targetColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $targetColumnWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $targetColumnWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $targetColumnWidth
